$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Calculated-column formulas used by Table1 for the Listening/Reading grade lookups
$lisFormula = '=IFERROR(INDEX(Sheet2!$F$5:$F$20, MATCH(Table1[[#This Row],[Lis_Mark]], Sheet2!$D$5:$D$20, 1)),"No Grade")'
$readFormula = '=IFERROR(INDEX(Sheet2!$F$5:$F$20, MATCH(Table1[[#This Row],[Read_Mark]], Sheet2!$D$5:$D$20, 1)),"No Grade")'

# Row 33 (Cambridge 13 Test 3 practice - Listening now scored)
$ws.Range("F33").Value = 34
$ws.Range("G33").Formula = $lisFormula
$ws.Range("H33").Value = 33
$ws.Range("I33").Formula = $readFormula
$ws.Range("K33").Value = 5

# Row 34
$ws.Range("F34").Value = 35
$ws.Range("G34").Formula = $lisFormula
$ws.Range("H34").Value = 30
$ws.Range("I34").Formula = $readFormula
$ws.Range("K34").Value = 5

# Row 35 (only Listening has been practiced so far; Read_Mark/Reading left blank)
$ws.Range("F35").Value = 31
$ws.Range("G35").Formula = $lisFormula
$ws.Range("K35").Value = 5

# Reflect the author's last active cell selection on the sheet
$ws.Range("J42").Select()
